# Generate Report for Handback
#
# - Status columns flip from "Ready for handoff" to
#   "Handed back: in sync with en-US" everywhere that text appears
#   (Overview!B2:C3 and the per-locale sheets' Status column).
# - Each locale sheet (zh-cn, de-de) gets its "Latest Target File" (F)
#   and "Latest Handback File" (G) columns populated with hyperlinks
#   mirroring the existing "Source File Name" (A) / "Latest Handoff
#   File" (D) links.
# - The "Latest Handback DateTime" (H) placeholder timestamp is
#   replaced with the real handback time for each locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status appears twice per row (zh-cn / de-de cols) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime,
        [string]$Row2TargetUrl,
        [string]$Row2HandbackUrl,
        [string]$Row3TargetUrl,
        [string]$Row3HandbackUrl
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Status -> Handed back
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Handback DateTime -> real timestamp instead of the
    # 0001-01-01 00:00:00 placeholder
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime

    # Latest Target File (F) / Latest Handback File (G), row 2
    # (mirrors A2 / D2 - the b831ba9f... file)
    $ws.Hyperlinks.Add($ws.Range("F2"), $Row2TargetUrl, [Type]::Missing, [Type]::Missing, $ws.Range("A2").Value())
    $ws.Hyperlinks.Add($ws.Range("G2"), $Row2HandbackUrl, [Type]::Missing, [Type]::Missing, $ws.Range("D2").Value())

    # Latest Target File (F) / Latest Handback File (G), row 3
    # (mirrors A3 / D3 - the f365c08c... file)
    $ws.Hyperlinks.Add($ws.Range("F3"), $Row3TargetUrl, [Type]::Missing, [Type]::Missing, $ws.Range("A3").Value())
    $ws.Hyperlinks.Add($ws.Range("G3"), $Row3HandbackUrl, [Type]::Missing, [Type]::Missing, $ws.Range("D3").Value())
}

Update-LocaleSheet "zh-cn" "2016-03-21 19:02:25" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a19b38252c42142d4bef91d170322f4458b09b54/e2e/b831ba9f-1757-464f-a51b-6d1ec7648a98.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b781e7a3b25ab459a410f23718678d4356e89a51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b831ba9f-1757-464f-a51b-6d1ec7648a98.5b12d99a3b0a97e63f4b57bebe90191e789ae057.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a19b38252c42142d4bef91d170322f4458b09b54/e2e/f365c08c-5968-4a31-b20e-6b0db5050c68.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b781e7a3b25ab459a410f23718678d4356e89a51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f365c08c-5968-4a31-b20e-6b0db5050c68.8aab68a2645f3777d27acf562308a1339066505f.zh-cn.xlf"

Update-LocaleSheet "de-de" "2016-03-21 19:02:31" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a19b38252c42142d4bef91d170322f4458b09b54/e2e/b831ba9f-1757-464f-a51b-6d1ec7648a98.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1818f284c0db9841dcd99bb91531fc8e9fbbd47e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b831ba9f-1757-464f-a51b-6d1ec7648a98.5b12d99a3b0a97e63f4b57bebe90191e789ae057.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a19b38252c42142d4bef91d170322f4458b09b54/e2e/f365c08c-5968-4a31-b20e-6b0db5050c68.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1818f284c0db9841dcd99bb91531fc8e9fbbd47e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f365c08c-5968-4a31-b20e-6b0db5050c68.8aab68a2645f3777d27acf562308a1339066505f.de-de.xlf"
